$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H119").Value = 1000
$ws.Range("J119").Value = 1000
$ws.Range("L119").Value = 3000
$ws.Range("N119").Value = -12676

$ws.Range("H121").Value = 3913.125
$ws.Range("J121").Value = 4157.857
$ws.Range("L121").Value = 12473.571
$ws.Range("N121").Value = -15967.571

$ws.Range("H132").Value = 2904.2
$ws.Range("I132").Value = 2847.6487
$ws.Range("J132").Value = 3601.6667
$ws.Range("K132").Value = 8542.946100000001
$ws.Range("L132").Value = 10805.0001
$ws.Range("M132").Value = -6012.946100000001
$ws.Range("N132").Value = -15865.0001

$ws.Range("H138").Value = 2339.4048
$ws.Range("J138").Value = 2078.9285
$ws.Range("L138").Value = 6236.7855
$ws.Range("N138").Value = -16516.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1187.7333
$ws.Range("I2").Value = 1234.75
$ws.Range("J2").Value = 999.6667
$ws.Range("K2").Value = 1234.75
$ws.Range("L2").Value = 999.6667
$ws.Range("M2").Value = -1121.75
$ws.Range("N2").Value = -1225.6667

$ws.Range("H32").Value = 860662.5600000001
$ws.Range("I32").Value = 1353101.5
$ws.Range("J32").Value = 16481.572
$ws.Range("K32").Value = 1353101.5
$ws.Range("L32").Value = 16481.572
$ws.Range("M32").Value = -1352814.5
$ws.Range("N32").Value = -17055.572

$ws.Range("H45").Value = 2777.625
$ws.Range("I45").Value = 2010.6428
$ws.Range("K45").Value = 2010.6428
$ws.Range("M45").Value = -1633.6428

$ws.Range("H61").Value = 2822.25
$ws.Range("I61").Value = 2080.6
$ws.Range("J61").Value = 3678
$ws.Range("K61").Value = 2080.6
$ws.Range("L61").Value = 3678
$ws.Range("M61").Value = -1868.6
$ws.Range("N61").Value = -4102

$ws.Range("H74").Value = 1786.7646
$ws.Range("I74").Value = 1696.8334
$ws.Range("J74").Value = 2002.6
$ws.Range("K74").Value = 1696.8334
$ws.Range("L74").Value = 2002.6
$ws.Range("M74").Value = -822.8334
$ws.Range("N74").Value = -3750.6

$ws.Range("H77").Value = 1786.7646
$ws.Range("I77").Value = 1696.8334
$ws.Range("J77").Value = 2002.6
$ws.Range("K77").Value = 8484.166999999999
$ws.Range("L77").Value = 10013
$ws.Range("M77").Value = -4116.166999999999
$ws.Range("N77").Value = -18749

$ws.Range("H110").Value = 2145
$ws.Range("I110").Value = 1931.25
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 1931.25
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = 113.75
$ws.Range("N110").Value = -7090

$ws.Range("H116").Value = 1187.7333
$ws.Range("I116").Value = 1234.75
$ws.Range("J116").Value = 999.6667
$ws.Range("K116").Value = 1234.75
$ws.Range("L116").Value = 999.6667
$ws.Range("M116").Value = 1059.25
$ws.Range("N116").Value = -5587.6667

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("N128").Value = $null
$ws.Range("L128").Value = 0

$ws.Range("H136").Value = 2822.25
$ws.Range("I136").Value = 2080.6
$ws.Range("J136").Value = 3678
$ws.Range("K136").Value = 6241.799999999999
$ws.Range("L136").Value = 11034
$ws.Range("M136").Value = -3691.799999999999
$ws.Range("N136").Value = -16134

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1187.7333
$ws.Range("I3").Value = 1234.75
$ws.Range("J3").Value = 999.6667
$ws.Range("K3").Value = 1234.75
$ws.Range("L3").Value = 999.6667
$ws.Range("M3").Value = -1120.75
$ws.Range("N3").Value = -1227.6667

$ws.Range("H134").Value = 3069
$ws.Range("I134").Value = 2885.7273
$ws.Range("K134").Value = 8657.1819
$ws.Range("M134").Value = -6122.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1421.2
$ws.Range("I58").Value = 783
$ws.Range("J58").Value = 2378.5
$ws.Range("K58").Value = 783
$ws.Range("L58").Value = 2378.5
$ws.Range("M58").Value = -580
$ws.Range("N58").Value = -2784.5

$ws.Range("H132").Value = 6174321
$ws.Range("I132").Value = 913.2857
$ws.Range("J132").Value = 27781248
$ws.Range("K132").Value = 2739.8571
$ws.Range("L132").Value = 83343744
$ws.Range("M132").Value = -209.8571000000002
$ws.Range("N132").Value = -83348804

$ws.Range("H134").Value = 2290.4443
$ws.Range("I134").Value = 2201.75
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 6605.25
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -4070.25
$ws.Range("N134").Value = -14070

$ws.Range("H136").Value = 1421.2
$ws.Range("I136").Value = 783
$ws.Range("J136").Value = 2378.5
$ws.Range("K136").Value = 2349
$ws.Range("L136").Value = 7135.5
$ws.Range("M136").Value = 201
$ws.Range("N136").Value = -12235.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 3316.6667
$ws.Range("I36").Value = 1300
$ws.Range("J36").Value = 5333.3335
$ws.Range("K36").Value = 3900
$ws.Range("L36").Value = 16000.0005
$ws.Range("M36").Value = -3731
$ws.Range("N36").Value = -16338.0005

$ws.Range("H60").Value = 2874.9822
$ws.Range("J60").Value = 3117.6274
$ws.Range("L60").Value = 9352.8822
$ws.Range("N60").Value = -9854.8822

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = $null

$ws.Range("H46").Value = 4349.95
$ws.Range("J46").Value = 4349.95
$ws.Range("L46").Value = 4349.95
$ws.Range("N46").Value = -4661.95

$ws.Range("H57").Value = 49530
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 49530
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = $null
$ws.Range("L57").Value = 49530
$ws.Range("N57").Value = -51170

$ws.Range("H80").Value = 31160574
$ws.Range("I80").Value = 57059052
$ws.Range("J80").Value = 82401.2
$ws.Range("K80").Value = 57059052
$ws.Range("L80").Value = 82401.2
$ws.Range("M80").Value = -57058054
$ws.Range("N80").Value = -84397.2

$ws.Range("H83").Value = 31160574
$ws.Range("I83").Value = 57059052
$ws.Range("J83").Value = 82401.2
$ws.Range("K83").Value = 285295260
$ws.Range("L83").Value = 412006
$ws.Range("M83").Value = -285290268
$ws.Range("N83").Value = -421990

$ws.Range("H126").Value = 3337
$ws.Range("I126").Value = 3112.6667
$ws.Range("J126").Value = 3673.5
$ws.Range("K126").Value = 9338.000100000001
$ws.Range("L126").Value = 11020.5
$ws.Range("M126").Value = -6868.000100000001
$ws.Range("N126").Value = -15960.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3192.875
$ws.Range("I136").Value = 4019.4
$ws.Range("J136").Value = 1815.3334
$ws.Range("K136").Value = 12058.2
$ws.Range("L136").Value = 5446.0002
$ws.Range("M136").Value = -9508.200000000001
$ws.Range("N136").Value = -10546.0002
